$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.21"
$ws.Range("E2").Value = "'1.66%"
$ws.Range("D3").Value = "'29.14"
$ws.Range("E3").Value = "'2.83%"
$ws.Range("D4").Value = "'5.207"
$ws.Range("E4").Value = "'3.11%"
$ws.Range("D5").Value = "'0.06980"
$ws.Range("E5").Value = "'7.08%"
$ws.Range("D6").Value = "'7.411"
$ws.Range("E6").Value = "'1.99%"
$ws.Range("D7").Value = "'3.555"
$ws.Range("E7").Value = "'5.59%"
$ws.Range("D8").Value = "'1.396"
$ws.Range("E8").Value = "'2.13%"
$ws.Range("D9").Value = "'0.8951"
$ws.Range("E9").Value = "'-3.69%"
$ws.Range("D10").Value = "'0.1616"
$ws.Range("E10").Value = "'4.11%"
$ws.Range("D11").Value = "'0.07501"
$ws.Range("E11").Value = "'22.73%"
$ws.Range("D12").Value = "'0.07772"
$ws.Range("E12").Value = "'2.34%"
$ws.Range("D13").Value = "'0.02932"
$ws.Range("E13").Value = "'1.16%"
$ws.Range("D14").Value = "'0.08992"
$ws.Range("E14").Value = "'0.15%"
$ws.Range("D15").Value = "'0.001590"
$ws.Range("E15").Value = "'0.14%"
$ws.Range("D16").Value = "'0.0006476"
$ws.Range("E16").Value = "'1.40%"
$ws.Range("D17").Value = "'0.006428"
$ws.Range("E17").Value = "'7.06%"
$ws.Range("D18").Value = "'3.461"
$ws.Range("E18").Value = "'0.27%"
$ws.Range("D19").Value = "'2.229"
$ws.Range("E19").Value = "'-0.29%"
$ws.Range("D20").Value = "'0.3234"
$ws.Range("E20").Value = "'0.95%"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("E21").Value = "'4.26%"
$ws.Range("E22").Value = "'-2.00%"
$ws.Range("D23").Value = "'0.1549"
$ws.Range("E23").Value = "'1.58%"
$ws.Range("D24").Value = "'0.04518"
$ws.Range("E24").Value = "'1.43%"
$ws.Range("D25").Value = "'0.001204"
$ws.Range("E25").Value = "'1.98%"
$ws.Range("D26").Value = "'0.004370"
$ws.Range("E26").Value = "'-0.52%"
$ws.Range("D27").Value = "'0.0001164"
$ws.Range("E27").Value = "'-6.72%"
$ws.Range("D28").Value = "'0.0001610"
$ws.Range("E28").Value = "'-0.30%"
$ws.Range("D40").Value = "'0.04340"
$ws.Range("E40").Value = "'4.47%"
$ws.Range("D41").Value = "'0.006928"
$ws.Range("E41").Value = "'4.32%"
$ws.Range("D42").Value = "'0.1247"
$ws.Range("E42").Value = "'1.90%"
$ws.Range("D43").Value = "'0.002069"
$ws.Range("E43").Value = "'2.62%"
$ws.Range("D44").Value = "'0.01171"
$ws.Range("E44").Value = "'-2.52%"
$ws.Range("D45").Value = "'0.00005839"
$ws.Range("E45").Value = "'5.86%"
$ws.Range("D47").Value = "'0.01301"
$ws.Range("E47").Value = "'0.23%"
